$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the job-specific data which gets reshuffled between
# rows 8, 9 and 10 (Activity, ScheduledResource, Plant, MISWorkCenter,
# PlannedResource).
$cols = "B", "D", "L", "O", "P"
$rows = @(8, 9, 10)

# Snapshot the current ("before") values for the three affected rows so the
# rotation below doesn't clobber data it still needs to read.
$before = @{}
foreach ($r in $rows) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $before[$addr] = $ws.Range($addr).Value2
    }
}

# Rotate the row data downward: row 9 gets row 8's old values, row 10 gets
# row 9's old values, and row 8 wraps around to row 10's old values.
$map = @{ 8 = 10; 9 = 8; 10 = 9 }

foreach ($r in $rows) {
    $src = $map[$r]
    foreach ($c in $cols) {
        $newValue = $before["$c$src"]
        $destAddr = "$c$r"
        # Skip no-op writes so cells whose content doesn't actually change
        # keep their original (e.g. empty shared-string) representation.
        if ($ws.Range($destAddr).Value2 -ne $newValue) {
            $ws.Range($destAddr).Value2 = $newValue
        }
    }
}

# Writing multi-line text (e.g. the PlannedResource column) can make the
# host auto-expand the row to a custom height. Auto-fit the touched rows
# back down so they return to the sheet's normal (non-custom) height.
foreach ($r in $rows) {
    $ws.Rows($r).AutoFit() | Out-Null
}
